$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-KeepRuns($para, $newText, $rPrXml) {
    # Rewrites only the paragraph's trailing run text (preserving any
    # leading empty run and the paragraph's own pPr) by inserting a fresh
    # <w:p> fragment into a Range that spans the old text but stops short
    # of the paragraph mark. This avoids the engine's normal behaviour of
    # silently coalescing a same-formatted empty run into the edited run
    # whenever a plain Find/Replace or Range.Text assignment touches it.
    $r = $para.Range
    $body = $d.Range($r.Start, $r.End - 1)
    if ($rPrXml) {
        $xml = "<w:p $wNs><w:r>$rPrXml<w:t>$newText</w:t></w:r></w:p>"
    } else {
        $xml = "<w:p $wNs><w:r><w:t>$newText</w:t></w:r></w:p>"
    }
    $body.InsertXML($xml)
}

# Map of the exact original paragraph text -> (new text, optional rPr XML).
# Using the full original paragraph text (rather than a global Find/Replace)
# keeps each edit scoped to exactly one paragraph, which matters here
# because "Play Great Rhino Megaways for Free: Review & RTP % 2021" occurs
# twice, in differently-formatted paragraphs.
$edits = @(
    @{ Old = "Play Great Rhino Megaways for Free: Review & RTP % 2021"; New = "Play Great Rhino Megaways for Free"; RPr = $null },
    @{ Old = "High maximum win per spin at 20,000x your bet"; New = "Cascade wins and increasing multipliers in free spins"; RPr = $null },
    @{ Old = "Up to 200,704 Megaways for increased chances of winning"; New = "High variance with potential for big wins"; RPr = $null },
    @{ Old = "Free spins with increasing multipliers and retrigger options"; New = "Engaging theme and graphics"; RPr = $null },
    @{ Old = "Autoplay feature and easy-to-use interface"; New = "User-friendly interface"; RPr = $null },
    @{ Old = "High volatility which may lead to dry spells"; New = "Limited bonus features"; RPr = $null },
    @{ Old = "Limited interactive bonus games"; New = "High variance may not appeal to all players"; RPr = $null },
    @{ Old = "Play Great Rhino Megaways for Free: Review & RTP % 2021"; New = "Play Great Rhino Megaways for Free"; RPr = "<w:rPr><w:b/></w:rPr>" },
    @{ Old = "Find out what we love about Great Rhino Megaways, a high-variance slot game with up to 200,704 Megaways and an RTP of over 96%. Play Great Rhino Megaways for free now."; New = "Read our review of Great Rhino Megaways and discover its gameplay features, theme, and potential winnings. Play for free!"; RPr = "<w:rPr><w:i/></w:rPr>" }
)

# Track which paragraphs have already been rewritten (by index) so that the
# duplicated heading text only consumes one matching paragraph per entry,
# in document order, even though two edits target the same old string.
$used = @{}

foreach ($edit in $edits) {
    $idx = 0
    foreach ($p in $d.Paragraphs) {
        if (-not $used.ContainsKey($idx)) {
            if ($p.Range.Text -eq ($edit.Old + "`r")) {
                Replace-KeepRuns $p $edit.New $edit.RPr
                $used[$idx] = $true
                break
            }
        }
        $idx = $idx + 1
    }
}
